$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Insert a new row at position 87 (this pushes former rows 87-95 down to
#    88-96, preserving their content/styles).
# ---------------------------------------------------------------------------
$ws.Rows.Item(87).Insert()

# ---------------------------------------------------------------------------
# 2) Populate the new row 87 ("Shared profile" monitoring rule).
#    The Insert() call copies some formatting down from row 86 onto the
#    blank row 87 (E87, J87, M87, O87, P87) - clean that up so the final
#    cell set/styles match the target exactly.
# ---------------------------------------------------------------------------

# A87 / B87 / C87 / E87 -> highlighted ("Y") style (fillId 4)
$ws.Range("A87").Value = "Shared profile"
$ws.Range("A87").Interior.Color = 3407820

$ws.Range("B87").Value = 1
$ws.Range("B87").Interior.Color = 3407820

$ws.Range("C87").Value = "Y"
$ws.Range("C87").Interior.Color = 3407820

$ws.Range("E87").Interior.Color = 3407820

# J87 - signal, no special style (drop the style inherited from the insert)
$ws.Range("J87").ClearFormats()
$ws.Range("J87").Value = "signal"

# K87 - "shared profile", Good style + left aligned
$ws.Range("K87").Value = "shared profile"
$ws.Range("K87").Style = "Good"
$ws.Range("K87").HorizontalAlignment = -4131

# O87 - blank, Good style + left aligned (re-apply explicitly, was inherited already)
$ws.Range("O87").Style = "Good"
$ws.Range("O87").HorizontalAlignment = -4131

# P87 - blank, Good style + wrap text (re-apply explicitly, was inherited already)
$ws.Range("P87").Style = "Good"
$ws.Range("P87").WrapText = $true

# M87 - inherited blank "Good" cell that should not exist in the final row
$ws.Range("M87").ClearFormats()
$ws.Range("M87").ClearContents()

# S87 - blank, Good style
$ws.Range("S87").Style = "Good"

# V87 - "shared profile", no special style
$ws.Range("V87").Value = "shared profile"

# ---------------------------------------------------------------------------
# 3) Row 86 gains an extra (blank, "Good"-styled) cell at S86.
# ---------------------------------------------------------------------------
$ws.Range("S86").Style = "Good"

# ---------------------------------------------------------------------------
# 4) The autofilter / filter-database range grows by the one inserted row
#    (from row 93 to row 94), while rows 95/96 (originally 94/95) stay
#    outside of it exactly as they did before the edit.
#    Re-applying AutoFilter() on this engine always snaps to the contiguous
#    "current region", so temporarily relocate the trailing rows out of the
#    way, set the filter on the now-exact range, then move them back.
# ---------------------------------------------------------------------------
$ws.Range("A95:Z96").Cut($ws.Range("A500:Z501"))

$ws.AutoFilterMode = $false
$ws.Range("A1:Z94").AutoFilter()

$ws.Range("A500:Z501").Cut($ws.Range("A95:Z96"))
$ws.Range("A500:Z501").ClearFormats()

# ---------------------------------------------------------------------------
# 5) Keep the workbook-level _FilterDatabase defined name in sync with the
#    new autofilter range.
# ---------------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Event params logged!_FilterDatabase") {
        $n.RefersTo = "='Event params logged'!`$A`$1:`$Z`$94"
    }
}
